$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 456.5
$ws.Range("J17").Value = 456.5
$ws.Range("L17").Value = 1369.5
$ws.Range("N17").Value = -1705.5
$ws.Range("H51").Value = 17106.334
$ws.Range("I51").Value = 15509
$ws.Range("K51").Value = 15509
$ws.Range("M51").Value = -15025
$ws.Range("H70").Value = 3249.375
$ws.Range("J70").Value = 3473.25
$ws.Range("L70").Value = 10419.75
$ws.Range("N70").Value = -10959.75
$ws.Range("H73").Value = 3249.375
$ws.Range("J73").Value = 3473.25
$ws.Range("L73").Value = 10419.75
$ws.Range("N73").Value = -12291.75
$ws.Range("H97").Value = 3356.7144
$ws.Range("J97").Value = 3249.5
$ws.Range("L97").Value = 9748.5
$ws.Range("N97").Value = -10740.5
$ws.Range("H111").Value = 2281.25
$ws.Range("I111").Value = 2182.375
$ws.Range("J111").Value = 2380.125
$ws.Range("K111").Value = 6547.125
$ws.Range("L111").Value = 7140.375
$ws.Range("M111").Value = -3480.125
$ws.Range("N111").Value = -13274.375
$ws.Range("H112").Value = 3504.5715
$ws.Range("I112").Value = 1316.3334
$ws.Range("J112").Value = 5145.75
$ws.Range("K112").Value = 3949.0002
$ws.Range("L112").Value = 15437.25
$ws.Range("M112").Value = -2841.0002
$ws.Range("N112").Value = -17653.25
$ws.Range("H132").Value = 5221.7617
$ws.Range("I132").Value = 4387.189
$ws.Range("K132").Value = 13161.567
$ws.Range("M132").Value = -10631.567
$ws.Range("H138").Value = 3473.024
$ws.Range("I138").Value = 6386.8335
$ws.Range("J138").Value = 2987.389
$ws.Range("K138").Value = 19160.5005
$ws.Range("L138").Value = 8962.167000000001
$ws.Range("M138").Value = -14020.5005
$ws.Range("N138").Value = -19242.167

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 2305.7837
$ws.Range("I74").Value = 2195.9355
$ws.Range("K74").Value = 2195.9355
$ws.Range("M74").Value = -1321.9355
$ws.Range("H77").Value = 2305.7837
$ws.Range("I77").Value = 2195.9355
$ws.Range("K77").Value = 10979.6775
$ws.Range("M77").Value = -6611.6775
$ws.Range("H122").Value = 2181.7932
$ws.Range("I122").Value = 1859.9375
$ws.Range("K122").Value = 5579.8125
$ws.Range("M122").Value = -3129.8125
$ws.Range("H132").Value = 40974.355
$ws.Range("I132").Value = 58877.777
$ws.Range("K132").Value = 176633.331
$ws.Range("M132").Value = -174103.331
$ws.Range("H135").Value = 54999.5
$ws.Range("J135").Value = 54999.5
$ws.Range("L135").Value = 54999.5
$ws.Range("N135").Value = -65139.5

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H59").Value = 79999
$ws.Range("J59").Value = 79999
$ws.Range("L59").Value = 79999
$ws.Range("N59").Value = -81693
$ws.Range("H94").Value = 2159.7144
$ws.Range("I94").Value = 1966
$ws.Range("K94").Value = 1966
$ws.Range("M94").Value = -1515
$ws.Range("H107").Value = 4400
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 4400
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4400
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -8240

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 4716
$ws.Range("I22").Value = 845.7778
$ws.Range("J22").Value = 13424
$ws.Range("K22").Value = 845.7778
$ws.Range("L22").Value = 13424
$ws.Range("M22").Value = -495.7778
$ws.Range("N22").Value = -14124
$ws.Range("H41").Value = 19045.455
$ws.Range("I41").Value = 3500
$ws.Range("J41").Value = 20600
$ws.Range("K41").Value = 3500
$ws.Range("L41").Value = 20600
$ws.Range("M41").Value = -3072
$ws.Range("N41").Value = -21456
$ws.Range("H132").Value = 1988
$ws.Range("I132").Value = 1988
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5964
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3434
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 71946.53
$ws.Range("I134").Value = 71946.53
$ws.Range("K134").Value = 215839.59
$ws.Range("M134").Value = -213304.59

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value = 1863849.1
$ws.Range("I4").Value = 1150234.1
$ws.Range("J4").Value = 9000000
$ws.Range("K4").Value = 3450702.3
$ws.Range("L4").Value = 27000000
$ws.Range("M4").Value = -3450590.3
$ws.Range("N4").Value = -27000224
$ws.Range("H12").Value = 135.09525
$ws.Range("J12").Value = 131.17647
$ws.Range("L12").Value = 393.52941
$ws.Range("N12").Value = -739.52941
$ws.Range("H28").Value = 1356.5
$ws.Range("I28").Value = 854.25
$ws.Range("K28").Value = 2562.75
$ws.Range("M28").Value = -2330.75
$ws.Range("H60").Value = 226.42857
$ws.Range("I60").Value = 249.16667
$ws.Range("J60").Value = 90
$ws.Range("K60").Value = 747.50001
$ws.Range("L60").Value = 270
$ws.Range("M60").Value = -496.50001
$ws.Range("N60").Value = -772
$ws.Range("H108").Value = 2689.4
$ws.Range("I108").Value = 1877.1111
$ws.Range("K108").Value = 5631.3333
$ws.Range("M108").Value = -2751.3333
$ws.Range("H114").Value = 994.6316
$ws.Range("I114").Value = 786.36365
$ws.Range("J114").Value = 1281
$ws.Range("K114").Value = 2359.09095
$ws.Range("L114").Value = 3843
$ws.Range("M114").Value = 894.9090500000002
$ws.Range("N114").Value = -10351
$ws.Range("H139").Value = 1228.1
$ws.Range("I139").Value = 1228.1
$ws.Range("K139").Value = 3684.3
$ws.Range("M139").Value = 1455.7

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H62").Value = 65000
$ws.Range("J62").Value = 65000
$ws.Range("L62").Value = 65000
$ws.Range("N62").Value = -66372
$ws.Range("H65").Value = 65000
$ws.Range("J65").Value = 65000
$ws.Range("L65").Value = 195000
$ws.Range("N65").Value = -201864
$ws.Range("H80").Value = 4929.8423
$ws.Range("J80").Value = 6484.5
$ws.Range("L80").Value = 6484.5
$ws.Range("N80").Value = -8480.5
$ws.Range("H83").Value = 4929.8423
$ws.Range("J83").Value = 6484.5
$ws.Range("L83").Value = 32422.5
$ws.Range("N83").Value = -42406.5
$ws.Range("H102").Value = 2809.7576
$ws.Range("I102").Value = 2590.3
$ws.Range("J102").Value = 5004.3335
$ws.Range("K102").Value = 2590.3
$ws.Range("L102").Value = 5004.3335
$ws.Range("M102").Value = -968.3000000000002
$ws.Range("N102").Value = -8248.333500000001
$ws.Range("H122").Value = 4312.25
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4312.25
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12936.75
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -17836.75
$ws.Range("H126").Value = 4363.593
$ws.Range("I126").Value = 4101.143
$ws.Range("J126").Value = 5282.1665
$ws.Range("K126").Value = 12303.429
$ws.Range("L126").Value = 15846.4995
$ws.Range("M126").Value = -9833.429
$ws.Range("N126").Value = -20786.4995

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 4935.143
$ws.Range("I40").Value = 4632.75
$ws.Range("K40").Value = 4632.75
$ws.Range("M40").Value = -4496.75
$ws.Range("H62").Value = 54500
$ws.Range("J62").Value = 54500
$ws.Range("L62").Value = 54500
$ws.Range("N62").Value = -55748
$ws.Range("H65").Value = 54500
$ws.Range("J65").Value = 54500
$ws.Range("L65").Value = 163500
$ws.Range("N65").Value = -169740
$ws.Range("H122").Value = 3977.25
$ws.Range("I122").Value = 3452
$ws.Range("J122").Value = 4502.5
$ws.Range("K122").Value = 10356
$ws.Range("L122").Value = 13507.5
$ws.Range("M122").Value = -7906
$ws.Range("N122").Value = -18407.5
$ws.Range("H133").Value = 85229.25
$ws.Range("J133").Value = 85229.25
$ws.Range("L133").Value = 85229.25
$ws.Range("N133").Value = -90289.25

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H46").Value = 50000
$ws.Range("J46").Value = 50000
$ws.Range("L46").Value = 50000
$ws.Range("N46").Value = -50462
$ws.Range("H48").Value = 24999.5
$ws.Range("J48").Value = 39999
$ws.Range("L48").Value = 39999
$ws.Range("N48").Value = -41137
$ws.Range("H49").Value = 44999.5
$ws.Range("J49").Value = 44999.5
$ws.Range("L49").Value = 44999.5
$ws.Range("N49").Value = -45459.5
$ws.Range("H63").Value = 39996.25
$ws.Range("J63").Value = 39996.25
$ws.Range("L63").Value = 39996.25
$ws.Range("N63").Value = -41244.25
$ws.Range("H66").Value = 39996.25
$ws.Range("J66").Value = 39996.25
$ws.Range("L66").Value = 119988.75
$ws.Range("N66").Value = -126228.75
$ws.Range("H132").Value = 58142.445
$ws.Range("I132").Value = 85488.914
$ws.Range("K132").Value = 256466.742
$ws.Range("M132").Value = -253936.742
$ws.Range("H133").Value = 81994.5
$ws.Range("J133").Value = 81994.5
$ws.Range("L133").Value = 81994.5
$ws.Range("N133").Value = -92114.5
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070
